$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 1551.2858
$ws.Cells.Item(43, 9).Value = 1653
$ws.Cells.Item(43, 11).Value = 1653
$ws.Cells.Item(43, 13).Value = -1584

$ws.Cells.Item(58, 8).Value = 1771.6666
$ws.Cells.Item(58, 9).Value = 907.5
$ws.Cells.Item(58, 11).Value = 2722.5
$ws.Cells.Item(58, 13).Value = -2572.5

$ws.Cells.Item(86, 8).Value = 1083
$ws.Cells.Item(86, 9).Value = 1023.2857
$ws.Cells.Item(86, 10).Value = 1501
$ws.Cells.Item(86, 11).Value = 1023.2857
$ws.Cells.Item(86, 12).Value = 1501
$ws.Cells.Item(86, 13).Value = 99.71429999999998
$ws.Cells.Item(86, 14).Value = -3747

$ws.Cells.Item(89, 8).Value = 1083
$ws.Cells.Item(89, 9).Value = 1023.2857
$ws.Cells.Item(89, 10).Value = 1501
$ws.Cells.Item(89, 11).Value = 5116.4285
$ws.Cells.Item(89, 12).Value = 7505
$ws.Cells.Item(89, 13).Value = 499.5715
$ws.Cells.Item(89, 14).Value = -18737

$ws.Cells.Item(121, 8).Value = 719.75
$ws.Cells.Item(121, 10).Value = 868.3333
$ws.Cells.Item(121, 12).Value = 2604.9999
$ws.Cells.Item(121, 14).Value = -6098.9999

$ws.Cells.Item(132, 8).Value = 834.95123
$ws.Cells.Item(132, 9).Value = 798.12823
$ws.Cells.Item(132, 10).Value = 1553
$ws.Cells.Item(132, 11).Value = 2394.38469
$ws.Cells.Item(132, 12).Value = 4659
$ws.Cells.Item(132, 13).Value = 135.6153100000001
$ws.Cells.Item(132, 14).Value = -9719

$ws.Cells.Item(138, 8).Value = 1625.9
$ws.Cells.Item(138, 9).Value = 1146.0968
$ws.Cells.Item(138, 10).Value = 1841.4637
$ws.Cells.Item(138, 11).Value = 3438.2904
$ws.Cells.Item(138, 12).Value = 5524.3911
$ws.Cells.Item(138, 13).Value = 1701.7096
$ws.Cells.Item(138, 14).Value = -15804.3911

$ws.Cells.Item(139, 8).Value = 69789.44500000001
$ws.Cells.Item(139, 10).Value = 69789.44500000001
$ws.Cells.Item(139, 12).Value = 69789.44500000001
$ws.Cells.Item(139, 14).Value = -80069.44500000001

$ws.Cells.Item(140, 8).Value = 84000
$ws.Cells.Item(140, 10).Value = 84000
$ws.Cells.Item(140, 12).Value = 84000
$ws.Cells.Item(140, 14).Value = -94360

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 264945.66
$ws.Cells.Item(2, 9).Value = 347546.8
$ws.Cells.Item(2, 11).Value = 347546.8
$ws.Cells.Item(2, 13).Value = -347433.8

$ws.Cells.Item(32, 8).Value = 1160.98
$ws.Cells.Item(32, 9).Value = 1160.98
$ws.Cells.Item(32, 11).Value = 1160.98
$ws.Cells.Item(32, 13).Value = -873.98

$ws.Cells.Item(74, 8).Value = 1282.7142
$ws.Cells.Item(74, 9).Value = 656.0909
$ws.Cells.Item(74, 10).Value = 3580.3333
$ws.Cells.Item(74, 11).Value = 656.0909
$ws.Cells.Item(74, 12).Value = 3580.3333
$ws.Cells.Item(74, 13).Value = 217.9091
$ws.Cells.Item(74, 14).Value = -5328.3333

$ws.Cells.Item(77, 8).Value = 1282.7142
$ws.Cells.Item(77, 9).Value = 656.0909
$ws.Cells.Item(77, 10).Value = 3580.3333
$ws.Cells.Item(77, 11).Value = 3280.4545
$ws.Cells.Item(77, 12).Value = 17901.6665
$ws.Cells.Item(77, 13).Value = 1087.5455
$ws.Cells.Item(77, 14).Value = -26637.6665

$ws.Cells.Item(102, 8).Value = 1608.3334
$ws.Cells.Item(102, 9).Value = 1608.3334
$ws.Cells.Item(102, 11).Value = 1608.3334
$ws.Cells.Item(102, 13).Value = 13.66660000000002

$ws.Cells.Item(116, 8).Value = 264945.66
$ws.Cells.Item(116, 9).Value = 347546.8
$ws.Cells.Item(116, 11).Value = 347546.8
$ws.Cells.Item(116, 13).Value = -345252.8

$ws.Cells.Item(122, 8).Value = 1697.3125
$ws.Cells.Item(122, 9).Value = 1154.0714
$ws.Cells.Item(122, 11).Value = 3462.2142
$ws.Cells.Item(122, 13).Value = -1012.2142

$ws.Cells.Item(132, 8).Value = 1852.6
$ws.Cells.Item(132, 9).Value = 1491
$ws.Cells.Item(132, 10).Value = 2603.6155
$ws.Cells.Item(132, 11).Value = 4473
$ws.Cells.Item(132, 12).Value = 7810.8465
$ws.Cells.Item(132, 13).Value = -1943
$ws.Cells.Item(132, 14).Value = -12870.8465

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 264945.66
$ws.Cells.Item(3, 9).Value = 347546.8
$ws.Cells.Item(3, 11).Value = 347546.8
$ws.Cells.Item(3, 13).Value = -347432.8

$ws.Cells.Item(20, 8).Value = 4368.8887
$ws.Cells.Item(20, 9).Value = 3831.75
$ws.Cells.Item(20, 10).Value = 4798.6
$ws.Cells.Item(20, 11).Value = 3831.75
$ws.Cells.Item(20, 12).Value = 4798.6
$ws.Cells.Item(20, 13).Value = -3584.75
$ws.Cells.Item(20, 14).Value = -5292.6

$ws.Cells.Item(99, 8).Value = 1551.6666
$ws.Cells.Item(99, 9).Value = 1670
$ws.Cells.Item(99, 11).Value = 1670
$ws.Cells.Item(99, 13).Value = -172

$ws.Cells.Item(107, 8).Value = 1311.375
$ws.Cells.Item(107, 9).Value = 1495.6
$ws.Cells.Item(107, 10).Value = 1004.3333
$ws.Cells.Item(107, 11).Value = 1495.6
$ws.Cells.Item(107, 12).Value = 1004.3333
$ws.Cells.Item(107, 13).Value = 424.4000000000001
$ws.Cells.Item(107, 14).Value = -4844.3333

$ws.Cells.Item(134, 8).Value = 6018.483
$ws.Cells.Item(134, 9).Value = 7439.905
$ws.Cells.Item(134, 10).Value = 2287.25
$ws.Cells.Item(134, 11).Value = 22319.715
$ws.Cells.Item(134, 12).Value = 6861.75
$ws.Cells.Item(134, 13).Value = -19784.715
$ws.Cells.Item(134, 14).Value = -11931.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 2240.7058
$ws.Cells.Item(86, 9).Value = 1980.5555
$ws.Cells.Item(86, 10).Value = 2533.375
$ws.Cells.Item(86, 11).Value = 1980.5555
$ws.Cells.Item(86, 12).Value = 2533.375
$ws.Cells.Item(86, 13).Value = -857.5554999999999
$ws.Cells.Item(86, 14).Value = -4779.375

$ws.Cells.Item(89, 8).Value = 2240.7058
$ws.Cells.Item(89, 9).Value = 1980.5555
$ws.Cells.Item(89, 10).Value = 2533.375
$ws.Cells.Item(89, 11).Value = 9902.7775
$ws.Cells.Item(89, 12).Value = 12666.875
$ws.Cells.Item(89, 13).Value = -4286.7775
$ws.Cells.Item(89, 14).Value = -23898.875

$ws.Cells.Item(107, 8).Value = 1049
$ws.Cells.Item(107, 10).Value = 275.75
$ws.Cells.Item(107, 12).Value = 275.75
$ws.Cells.Item(107, 14).Value = -4115.75

$ws.Cells.Item(134, 8).Value = 2450.4348
$ws.Cells.Item(134, 9).Value = 2352.5
$ws.Cells.Item(134, 11).Value = 7057.5
$ws.Cells.Item(134, 13).Value = -4522.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(141, 8).Value = 3596.25
$ws.Cells.Item(141, 9).Value = 3500.4666
$ws.Cells.Item(141, 10).Value = 5033
$ws.Cells.Item(141, 11).Value = 10501.3998
$ws.Cells.Item(141, 12).Value = 15099
$ws.Cells.Item(141, 13).Value = -5321.399800000001
$ws.Cells.Item(141, 14).Value = -25459

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4153.1816
$ws.Cells.Item(70, 9).Value = 3909.5
$ws.Cells.Item(70, 10).Value = 4579.625
$ws.Cells.Item(70, 11).Value = 3909.5
$ws.Cells.Item(70, 12).Value = 4579.625
$ws.Cells.Item(70, 13).Value = -3639.5
$ws.Cells.Item(70, 14).Value = -5119.625

$ws.Cells.Item(73, 8).Value = 4153.1816
$ws.Cells.Item(73, 9).Value = 3909.5
$ws.Cells.Item(73, 10).Value = 4579.625
$ws.Cells.Item(73, 11).Value = 3909.5
$ws.Cells.Item(73, 12).Value = 4579.625
$ws.Cells.Item(73, 13).Value = -2973.5
$ws.Cells.Item(73, 14).Value = -6451.625

$ws.Cells.Item(97, 8).Value = 957.75757
$ws.Cells.Item(97, 9).Value = 946.9231
$ws.Cells.Item(97, 11).Value = 946.9231
$ws.Cells.Item(97, 13).Value = -450.9231

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3527.1667
$ws.Cells.Item(7, 9).Value = 3590.6
$ws.Cells.Item(7, 11).Value = 3590.6
$ws.Cells.Item(7, 13).Value = -3478.6

$ws.Cells.Item(16, 8).Value = 5340.4614
$ws.Cells.Item(16, 9).Value = 5748
$ws.Cells.Item(16, 11).Value = 5748
$ws.Cells.Item(16, 13).Value = -5578

$ws.Cells.Item(93, 8).Value = 919.5
$ws.Cells.Item(93, 9).Value = 503.24
$ws.Cells.Item(93, 11).Value = 503.24
$ws.Cells.Item(93, 13).Value = 744.76

$ws.Cells.Item(126, 8).Value = 3527.1667
$ws.Cells.Item(126, 9).Value = 3590.6
$ws.Cells.Item(126, 11).Value = 10771.8
$ws.Cells.Item(126, 13).Value = -8301.799999999999

$ws.Cells.Item(136, 8).Value = 3993.7273
$ws.Cells.Item(136, 9).Value = 3319.5
$ws.Cells.Item(136, 10).Value = 4802.8
$ws.Cells.Item(136, 11).Value = 9958.5
$ws.Cells.Item(136, 12).Value = 14408.4
$ws.Cells.Item(136, 13).Value = -7408.5
$ws.Cells.Item(136, 14).Value = -19508.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 5036.3335
$ws.Cells.Item(81, 9).Value = 1709.4
$ws.Cells.Item(81, 10).Value = 6699.8
$ws.Cells.Item(81, 11).Value = 3418.8
$ws.Cells.Item(81, 12).Value = 13399.6
$ws.Cells.Item(81, 13).Value = -2357.8
$ws.Cells.Item(81, 14).Value = -15521.6

$ws.Cells.Item(84, 8).Value = 5036.3335
$ws.Cells.Item(84, 9).Value = 1709.4
$ws.Cells.Item(84, 10).Value = 6699.8
$ws.Cells.Item(84, 11).Value = 17094
$ws.Cells.Item(84, 12).Value = 66998
$ws.Cells.Item(84, 13).Value = -11790
$ws.Cells.Item(84, 14).Value = -77606

$ws.Cells.Item(113, 8).Value = 574.75
$ws.Cells.Item(113, 9).Value = 421
$ws.Cells.Item(113, 10).Value = 728.5
$ws.Cells.Item(113, 11).Value = 1263
$ws.Cells.Item(113, 12).Value = 2185.5
$ws.Cells.Item(113, 13).Value = 907
$ws.Cells.Item(113, 14).Value = -6525.5

$ws.Cells.Item(126, 8).Value = 2313.9443
$ws.Cells.Item(126, 9).Value = 2215.7
$ws.Cells.Item(126, 10).Value = 2436.75
$ws.Cells.Item(126, 11).Value = 6647.099999999999
$ws.Cells.Item(126, 12).Value = 7310.25
$ws.Cells.Item(126, 13).Value = -4177.099999999999
$ws.Cells.Item(126, 14).Value = -12250.25

$ws.Cells.Item(132, 8).Value = 1154.34
$ws.Cells.Item(132, 9).Value = 1015.04877
$ws.Cells.Item(132, 10).Value = 1788.8889
$ws.Cells.Item(132, 11).Value = 3045.14631
$ws.Cells.Item(132, 12).Value = 5366.6667
$ws.Cells.Item(132, 13).Value = -515.1463100000001
$ws.Cells.Item(132, 14).Value = -10426.6667
